# The grid-point location (and the T2 temperature series sampled at it)
# changed. Convert every old value to its replacement with a literal,
# match-case, whole-document Find/Replace — each old number/string is
# unique in the doc, so a straight Execute(..., Replace:=2) per pair is
# enough; no wildcards needed.
$d = $word.ActiveDocument

$d.Content.Find.Execute("Grid Point: Lat 1.0905, Lon 36.4858", $true, $false, $false, $false, $false, $true, 1, $false, "Grid Point: Lat 2.5989, Lon 40.1917", 2) | Out-Null
$d.Content.Find.Execute("Dropped constant columns: lon(36.48581632416828), lat(1.090492598507372)", $true, $false, $false, $false, $false, $true, 1, $false, "Dropped constant columns: lon(40.19166022991316), lat(2.5989427479035143)", 2) | Out-Null
$d.Content.Find.Execute("-269.66544", $true, $false, $false, $false, $false, $true, 1, $false, "2.1087952", 2) | Out-Null
$d.Content.Find.Execute("-270.3573", $true, $false, $false, $false, $false, $true, 1, $false, "3.3631897", 2) | Out-Null
$d.Content.Find.Execute("-270.52338", $true, $false, $false, $false, $false, $true, 1, $false, "3.4154053", 2) | Out-Null
$d.Content.Find.Execute("-272.6171", $true, $false, $false, $false, $false, $true, 1, $false, "3.3013916", 2) | Out-Null
$d.Content.Find.Execute("-271.13763", $true, $false, $false, $false, $false, $true, 1, $false, "3.3478699", 2) | Out-Null
$d.Content.Find.Execute("-271.79184", $true, $false, $false, $false, $false, $true, 1, $false, "3.3709412", 2) | Out-Null
$d.Content.Find.Execute("-272.23157", $true, $false, $false, $false, $false, $true, 1, $false, "3.0004578", 2) | Out-Null
$d.Content.Find.Execute("-273.16678", $true, $false, $false, $false, $false, $true, 1, $false, "2.589203", 2) | Out-Null
$d.Content.Find.Execute("-275.66473", $true, $false, $false, $false, $false, $true, 1, $false, "1.6613464", 2) | Out-Null
$d.Content.Find.Execute("-273.7265", $true, $false, $false, $false, $false, $true, 1, $false, "0.43762207", 2) | Out-Null
$d.Content.Find.Execute("-275.51608", $true, $false, $false, $false, $false, $true, 1, $false, "-1.5089417", 2) | Out-Null
$d.Content.Find.Execute("-274.16443", $true, $false, $false, $false, $false, $true, 1, $false, "-2.698639", 2) | Out-Null
$d.Content.Find.Execute("-276.18225", $true, $false, $false, $false, $false, $true, 1, $false, "-3.005371", 2) | Out-Null
$d.Content.Find.Execute("-275.0824", $true, $false, $false, $false, $false, $true, 1, $false, "-2.8447876", 2) | Out-Null
$d.Content.Find.Execute("-275.00455", $true, $false, $false, $false, $false, $true, 1, $false, "-2.6717224", 2) | Out-Null
$d.Content.Find.Execute("-274.65015", $true, $false, $false, $false, $false, $true, 1, $false, "-2.804596", 2) | Out-Null
$d.Content.Find.Execute("-274.39697", $true, $false, $false, $false, $false, $true, 1, $false, "-2.6209717", 2) | Out-Null
$d.Content.Find.Execute("-274.2068", $true, $false, $false, $false, $false, $true, 1, $false, "-2.1082458", 2) | Out-Null
$d.Content.Find.Execute("-274.03665", $true, $false, $false, $false, $false, $true, 1, $false, "-1.1530151", 2) | Out-Null
$d.Content.Find.Execute("-273.78427", $true, $false, $false, $false, $false, $true, 1, $false, "-1.9770203", 2) | Out-Null
$d.Content.Find.Execute("-273.48276", $true, $false, $false, $false, $false, $true, 1, $false, "-2.23172", 2) | Out-Null
$d.Content.Find.Execute("-272.77905", $true, $false, $false, $false, $false, $true, 1, $false, "-3.1305847", 2) | Out-Null
$d.Content.Find.Execute("-272.31378", $true, $false, $false, $false, $false, $true, 1, $false, "-1.3967285", 2) | Out-Null
$d.Content.Find.Execute("-269.39874", $true, $false, $false, $false, $false, $true, 1, $false, "1.036499", 2) | Out-Null
$d.Content.Find.Execute("-267.69788", $true, $false, $false, $false, $false, $true, 1, $false, "4.049164", 2) | Out-Null
$d.Content.Find.Execute("-267.39346", $true, $false, $false, $false, $false, $true, 1, $false, "5.806885", 2) | Out-Null
$d.Content.Find.Execute("-270.50082", $true, $false, $false, $false, $false, $true, 1, $false, "5.5571594", 2) | Out-Null
$d.Content.Find.Execute("-272.3076", $true, $false, $false, $false, $false, $true, 1, $false, "4.6316833", 2) | Out-Null
$d.Content.Find.Execute("-272.05832", $true, $false, $false, $false, $false, $true, 1, $false, "3.8804932", 2) | Out-Null
$d.Content.Find.Execute("-270.9744", $true, $false, $false, $false, $false, $true, 1, $false, "3.2436523", 2) | Out-Null
$d.Content.Find.Execute("-270.71457", $true, $false, $false, $false, $false, $true, 1, $false, "2.671753", 2) | Out-Null
$d.Content.Find.Execute("-271.3162", $true, $false, $false, $false, $false, $true, 1, $false, "2.1029663", 2) | Out-Null
$d.Content.Find.Execute("-272.2075", $true, $false, $false, $false, $false, $true, 1, $false, "1.3392944", 2) | Out-Null
$d.Content.Find.Execute("-273.57837", $true, $false, $false, $false, $false, $true, 1, $false, "-0.26794434", 2) | Out-Null
$d.Content.Find.Execute("-276.08368", $true, $false, $false, $false, $false, $true, 1, $false, "-2.3234863", 2) | Out-Null
$d.Content.Find.Execute("-277.19434", $true, $false, $false, $false, $false, $true, 1, $false, "-3.52063", 2) | Out-Null
$d.Content.Find.Execute("-277.24664", $true, $false, $false, $false, $false, $true, 1, $false, "-3.6307678", 2) | Out-Null
$d.Content.Find.Execute("-275.85336", $true, $false, $false, $false, $false, $true, 1, $false, "-2.9499207", 2) | Out-Null
$d.Content.Find.Execute("-275.07416", $true, $false, $false, $false, $false, $true, 1, $false, "-3.3542786", 2) | Out-Null
$d.Content.Find.Execute("-274.5754", $true, $false, $false, $false, $false, $true, 1, $false, "-3.3404236", 2) | Out-Null
$d.Content.Find.Execute("-274.06244", $true, $false, $false, $false, $false, $true, 1, $false, "-2.542633", 2) | Out-Null
$d.Content.Find.Execute("-274.2919", $true, $false, $false, $false, $false, $true, 1, $false, "-1.3841248", 2) | Out-Null
$d.Content.Find.Execute("-274.68628", $true, $false, $false, $false, $false, $true, 1, $false, "-0.6476135", 2) | Out-Null
$d.Content.Find.Execute("-274.90677", $true, $false, $false, $false, $false, $true, 1, $false, "-1.7875977", 2) | Out-Null
$d.Content.Find.Execute("-274.6269", $true, $false, $false, $false, $false, $true, 1, $false, "-2.4434814", 2) | Out-Null
$d.Content.Find.Execute("-273.68192", $true, $false, $false, $false, $false, $true, 1, $false, "-2.8603516", 2) | Out-Null
$d.Content.Find.Execute("-273.06653", $true, $false, $false, $false, $false, $true, 1, $false, "-1.6043701", 2) | Out-Null
$d.Content.Find.Execute("-269.88913", $true, $false, $false, $false, $false, $true, 1, $false, "0.8451233", 2) | Out-Null
$d.Content.Find.Execute("-267.49268", $true, $false, $false, $false, $false, $true, 1, $false, "3.479187", 2) | Out-Null
$d.Content.Find.Execute("-267.387", $true, $false, $false, $false, $false, $true, 1, $false, "5.3216553", 2) | Out-Null
$d.Content.Find.Execute("-270.1719", $true, $false, $false, $false, $false, $true, 1, $false, "5.441498", 2) | Out-Null
$d.Content.Find.Execute("-272.0983", $true, $false, $false, $false, $false, $true, 1, $false, "5.0298157", 2) | Out-Null
$d.Content.Find.Execute("-271.4964", $true, $false, $false, $false, $false, $true, 1, $false, "4.4765015", 2) | Out-Null
$d.Content.Find.Execute("-270.902", $true, $false, $false, $false, $false, $true, 1, $false, "3.960144", 2) | Out-Null
$d.Content.Find.Execute("-270.97256", $true, $false, $false, $false, $false, $true, 1, $false, "2.9674377", 2) | Out-Null
$d.Content.Find.Execute("-271.54843", $true, $false, $false, $false, $false, $true, 1, $false, "2.1362915", 2) | Out-Null
$d.Content.Find.Execute("-272.70773", $true, $false, $false, $false, $false, $true, 1, $false, "0.9539795", 2) | Out-Null
$d.Content.Find.Execute("-274.75095", $true, $false, $false, $false, $false, $true, 1, $false, "-0.45114136", 2) | Out-Null
$d.Content.Find.Execute("-277.0655", $true, $false, $false, $false, $false, $true, 1, $false, "-2.5144043", 2) | Out-Null
$d.Content.Find.Execute("-277.87247", $true, $false, $false, $false, $false, $true, 1, $false, "-3.4780273", 2) | Out-Null
$d.Content.Find.Execute("-277.41455", $true, $false, $false, $false, $false, $true, 1, $false, "-3.6743774", 2) | Out-Null
$d.Content.Find.Execute("-276.07465", $true, $false, $false, $false, $false, $true, 1, $false, "-3.2113647", 2) | Out-Null
$d.Content.Find.Execute("-275.50256", $true, $false, $false, $false, $false, $true, 1, $false, "-2.7855835", 2) | Out-Null
$d.Content.Find.Execute("-275.14078", $true, $false, $false, $false, $false, $true, 1, $false, "-2.3372192", 2) | Out-Null
$d.Content.Find.Execute("-274.8572", $true, $false, $false, $false, $false, $true, 1, $false, "-2.1173706", 2) | Out-Null
$d.Content.Find.Execute("-274.5533", $true, $false, $false, $false, $false, $true, 1, $false, "-3.0126038", 2) | Out-Null
$d.Content.Find.Execute("-274.20456", $true, $false, $false, $false, $false, $true, 1, $false, "-3.270935", 2) | Out-Null
$d.Content.Find.Execute("-273.8013", $true, $false, $false, $false, $false, $true, 1, $false, "-3.1472778", 2) | Out-Null
$d.Content.Find.Execute("-273.343", $true, $false, $false, $false, $false, $true, 1, $false, "-2.3453064", 2) | Out-Null
$d.Content.Find.Execute("-272.85403", $true, $false, $false, $false, $false, $true, 1, $false, "-2.1729126", 2) | Out-Null
$d.Content.Find.Execute("-272.6955", $true, $false, $false, $false, $false, $true, 1, $false, "-1.5320129", 2) | Out-Null
$d.Content.Find.Execute("-270.31418", $true, $false, $false, $false, $false, $true, 1, $false, "0.40621948", 2) | Out-Null
$d.Content.Find.Execute("-268.71185", $true, $false, $false, $false, $false, $true, 1, $false, "2.765686", 2) | Out-Null
